# Export model responses for checking, and filter eval results by valid
# questions only: add a "valid" helper column (S) that flags rows where
# all four check columns (F, G, H, I) are "yes".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("S1").Value = "valid"

# Formula for each data row (2..181) -- rows below 181 have no F:I data
$ws.Range("S2:S181").Formula = '=INT(AND(F2="yes", G2="yes", H2="yes",I2="yes"))'
